$wb = $excel.ActiveWorkbook

# --- Update "Schedule" sheet ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 703.1615954999999
$schedule.Range("F2").Value = 11.62634913194444

# --- Update "Detailed" sheet ---
$detailed = $wb.Worksheets.Item("Detailed")

$updates = @{
    20 = @{ B = 0 }
    21 = @{ B = -0.9298; C = "historical" }
    22 = @{ B = -5.45024; C = "historical" }
    23 = @{ B = -5.58973 }
    24 = @{ B = -5.01 }
    25 = @{ B = -5.70675 }
    26 = @{ B = -0.94965 }
    27 = @{ B = -0.9498799999999999 }
    28 = @{ B = -2.84053 }
    29 = @{ B = -0.99452 }
    30 = @{ B = 0.00027 }
    31 = @{ B = 0.51 }
    32 = @{ B = 2.43896 }
    33 = @{ B = 8.276960000000001 }
    34 = @{ B = 8.31124 }
    35 = @{ B = 0 }
    36 = @{ B = 0.264 }
    37 = @{ B = 4.52793 }
    38 = @{ B = 25.06795 }
    39 = @{ B = 43.32482 }
    40 = @{ B = 58.39741 }
    41 = @{ B = 62.41913 }
    42 = @{ B = 59.18831 }
    43 = @{ B = 65 }
    44 = @{ B = 69.22541 }
    45 = @{ B = 65 }
    46 = @{ B = 65 }
    47 = @{ B = 65 }
    48 = @{ B = 64.43608999999999 }
    49 = @{ B = 64.01244 }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    if ($rowData.ContainsKey("B")) {
        $detailed.Range("B$row").Value = $rowData["B"]
    }
    if ($rowData.ContainsKey("C")) {
        $detailed.Range("C$row").Value = $rowData["C"]
    }
}
